# Updated symbol list on Thu Dec 29 21:20:23 UTC 2022 with GitHub Actions
# Applies the per-row "Price" (column D) updates and the two "Volume(1h)"
# (column E) label corrections described by the source diff.
#
# All D-column cells in this sheet are stored as literal text (inlineStr),
# e.g. "245.74", "0.001525" - not numbers. Plain `Range.Value = "245.30"`
# would be auto-coerced to a number by Excel (losing the trailing zero /
# changing the stored type), so we force the cell's number format to Text
# ("@") before writing the literal string, which keeps it a text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# Column D ("Price") updates
Set-TextValue "D2"  "245.30"
Set-TextValue "D3"  "24.19"
Set-TextValue "D4"  "5.281"
Set-TextValue "D7"  "3.147"
Set-TextValue "D8"  "0.8159"
Set-TextValue "D9"  "0.8449"
Set-TextValue "D10" "0.1362"
Set-TextValue "D11" "0.06960"
Set-TextValue "D12" "0.03129"
Set-TextValue "D13" "0.02900"
Set-TextValue "D14" "0.09386"
Set-TextValue "D16" "0.001524"
Set-TextValue "D18" "0.0005981"
Set-TextValue "D19" "0.006094"
Set-TextValue "D21" "0.004614"
Set-TextValue "D22" "0.00006901"
Set-TextValue "D23" "3.501"
Set-TextValue "D24" "2.147"
Set-TextValue "D27" "0.1359"
Set-TextValue "D28" "0.0002332"
Set-TextValue "D41" "0.006281"
Set-TextValue "D42" "0.1052"
Set-TextValue "D43" "0.003401"
Set-TextValue "D44" "0.008500"
Set-TextValue "D47" "0.3701"
Set-TextValue "D48" "0.002288"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"

# Column E ("Volume(1h)") label corrections
Set-TextValue "E21" "20HotbitTokenHTB"
Set-TextValue "E41" "40KickTokenKICKBestin24h"
